$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    for ($row = 2; $row -le 35; $row++) {
        $cell = $ws.Cells.Item($row, 4)
        if ($cell.Text -eq "(0, 0)") {
            $cell.Value = "(nan, nan)"
        }
    }
}
